$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the greeting text for the "R10" rule row
$ws.Range("E8").Value = "GIT UPDATE"

# Select the edited cell, matching the saved selection state
$ws.Range("E8").Select()
